# Apply the crypto-price/volume refresh described by the commit diff.
# Source data is plain text (coin name / link / price / % volume) stored as
# inline strings, so every assignment below targets Range.Value directly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" cells (column D) are plain decimals (e.g. "0.06448") that Excel
# would otherwise auto-convert to a Number on assignment. Force those specific
# cells to Text first, then restore the default "Normal" style afterwards so the
# cell keeps matching its original (unstyled) formatting.

$ws.Range("D2").Value = "26.221.07"
$ws.Range("E2").Value = "  -3.73%  "

$ws.Range("D3").Value = "1.660.32"
$ws.Range("E3").Value = "  -2.48%  "

$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.45%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5152"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.82%  "

$ws.Range("E7").Value = "  +0.35%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2580"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.86%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06448"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.59%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07815"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.56%  "

$ws.Range("D12").Value = "1.662.33"
$ws.Range("E12").Value = "  -2.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.300"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.29%  "

$ws.Range("D14").Value = "1.887.72"
$ws.Range("E14").Value = "  -2.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5549"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.16%  "

$ws.Range("D16").Value = "0.0₅8069"
$ws.Range("E16").Value = "  -0.94%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.48%  "

$ws.Range("D18").Value = "26.246.79"
$ws.Range("E18").Value = "  -3.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "212.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.37%  "

$ws.Range("E20").Value = "  +0.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.433"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.82%  "

$ws.Range("E22").Value = "  -2.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.022"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.91%  "

$ws.Range("E24").Value = "  +0.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.756"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.23%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1173"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.82%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.995"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.81%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05214"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.254"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.361"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.89%  "

$ws.Range("E33").Value = "  -4.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.576"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.76%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9344"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.98%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.372"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.36%  "

$ws.Range("D38").Value = "1.175.02"
$ws.Range("E38").Value = "  +13.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5699"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.00%  "

$ws.Range("E40").Value = "  -2.26%  "

$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.005"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.29%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8408"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.686"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.18%  "

$ws.Range("D45").Value = "1.797.79"
$ws.Range("E45").Value = "  -2.61%  "

$ws.Range("D46").Value = "0.0₈115"
$ws.Range("E46").Value = "  +0.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4539"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "56.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.79%  "

$ws.Range("E49").Value = "  -0.31%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.902"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05060"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.18%  "
